$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the header (row 1), pushing existing data down.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new LEM/LED pair.
$ws.Cells.Item(2, 1).Value = "LEM-198-16-3022KH"
$ws.Cells.Item(2, 2).Value = "LED-198-H35-3022"

# Re-sort the data (below the header in row 1) by column A, ascending, to keep alphabetical order.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A1:A49"))
$sortObj.SetRange($ws.Range("A1:B49"))
$sortObj.Header = 1
$sortObj.Apply()

# Update selection to match the saved view state.
$ws.Range("B3").Select()
